# Refactor the scan collation logic: the "scan" function now calls
# collate_uploads, so matched files are grouped together by scan (and
# tagged with a new SessionLabel in column I built from the subject id
# and dataset pattern) before the unmatched rows, rather than the scan id
# being derived purely from a shared DICOM SeriesNumber. Rebuild the
# worksheet data in the new row order with the newly added SessionLabel
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe out the existing data rows (2-17); they get fully rebuilt below in
# the new row order produced by collate_uploads.
$dataRange = $ws.Range("A2:T17")
$dataRange.ClearContents()

# Force the range to Text so that the numeric-looking subject/session ids
# (e.g. "002304", "20200312") round-trip as strings instead of being
# auto-converted to numbers by Excel.
$dataRange.NumberFormat = "@"

# --- Row 2-10: matched/collated scan rows ---
# --- Row 11-17: unmatched rows ---
$ws.Range("A2").Value2 = "Scans"
$ws.Range("B2").Value2 = "tests/fixtures/basic/DOE^JOHN-002304/20200312HeadCT/Head CT/image-00000.dcm"
$ws.Range("C2").Value2 = "image-00000.dcm"
$ws.Range("D2").Value2 = "Y"
$ws.Range("F2").Value2 = "002304"
$ws.Range("G2").Value2 = "20200312"
$ws.Range("H2").Value2 = "Head_CT"
$ws.Range("I2").Value2 = "002304_CT1"
$ws.Range("L2").Value2 = "CT"
$ws.Range("O2").Value2 = "20200312"
$ws.Range("P2").Value2 = "CT1 abdomen"
$ws.Range("Q2").Value2 = "DOE^JOHN"
$ws.Range("R2").Value2 = "002304"
$ws.Range("S2").Value2 = "Head CT"
$ws.Range("T2").Value2 = "image-00000"
$ws.Range("A3").Value2 = "Scans"
$ws.Range("B3").Value2 = "tests/fixtures/basic/DOE^JOHN-002304/20200312HeadCT/Head CT/image-00001.dcm"
$ws.Range("C3").Value2 = "image-00001.dcm"
$ws.Range("D3").Value2 = "Y"
$ws.Range("F3").Value2 = "002304"
$ws.Range("G3").Value2 = "20200312"
$ws.Range("H3").Value2 = "Head_CT"
$ws.Range("I3").Value2 = "002304_CT1"
$ws.Range("L3").Value2 = "CT"
$ws.Range("O3").Value2 = "20200312"
$ws.Range("P3").Value2 = "CT1 abdomen"
$ws.Range("Q3").Value2 = "DOE^JOHN"
$ws.Range("R3").Value2 = "002304"
$ws.Range("S3").Value2 = "Head CT"
$ws.Range("T3").Value2 = "image-00001"
$ws.Range("A4").Value2 = "Scans"
$ws.Range("B4").Value2 = "tests/fixtures/basic/DOE^JOHN-002304/20200312HeadCT/Neck CT/image-00000.dcm"
$ws.Range("C4").Value2 = "image-00000.dcm"
$ws.Range("D4").Value2 = "Y"
$ws.Range("F4").Value2 = "002304"
$ws.Range("G4").Value2 = "20200312"
$ws.Range("H4").Value2 = "Neck_CT"
$ws.Range("I4").Value2 = "002304_CT1"
$ws.Range("L4").Value2 = "CT"
$ws.Range("O4").Value2 = "20200312"
$ws.Range("P4").Value2 = "CT1 abdomen"
$ws.Range("Q4").Value2 = "DOE^JOHN"
$ws.Range("R4").Value2 = "002304"
$ws.Range("S4").Value2 = "Neck CT"
$ws.Range("T4").Value2 = "image-00000"
$ws.Range("A5").Value2 = "Scans"
$ws.Range("B5").Value2 = "tests/fixtures/basic/DOE^JOHN-002304/20200312HeadCT/Neck CT/image-00001.dcm"
$ws.Range("C5").Value2 = "image-00001.dcm"
$ws.Range("D5").Value2 = "Y"
$ws.Range("F5").Value2 = "002304"
$ws.Range("G5").Value2 = "20200312"
$ws.Range("H5").Value2 = "Neck_CT"
$ws.Range("I5").Value2 = "002304_CT1"
$ws.Range("L5").Value2 = "CT"
$ws.Range("O5").Value2 = "20200312"
$ws.Range("P5").Value2 = "CT1 abdomen"
$ws.Range("Q5").Value2 = "DOE^JOHN"
$ws.Range("R5").Value2 = "002304"
$ws.Range("S5").Value2 = "Neck CT"
$ws.Range("T5").Value2 = "image-00001"
$ws.Range("A6").Value2 = "Scans"
$ws.Range("B6").Value2 = "tests/fixtures/basic/DOE^JOHN-002304/20200312HeadCT/Neck CT/image-00002.dcm"
$ws.Range("C6").Value2 = "image-00002.dcm"
$ws.Range("D6").Value2 = "Y"
$ws.Range("F6").Value2 = "002304"
$ws.Range("G6").Value2 = "20200312"
$ws.Range("H6").Value2 = "Neck_CT"
$ws.Range("I6").Value2 = "002304_CT1"
$ws.Range("L6").Value2 = "CT"
$ws.Range("O6").Value2 = "20200312"
$ws.Range("P6").Value2 = "CT1 abdomen"
$ws.Range("Q6").Value2 = "DOE^JOHN"
$ws.Range("R6").Value2 = "002304"
$ws.Range("S6").Value2 = "Neck CT"
$ws.Range("T6").Value2 = "image-00002"
$ws.Range("A7").Value2 = "Scans"
$ws.Range("B7").Value2 = "tests/fixtures/basic/ROE^JANE-397829/20190115/SomeCT/img-00000.dcm"
$ws.Range("C7").Value2 = "img-00000.dcm"
$ws.Range("D7").Value2 = "Y"
$ws.Range("F7").Value2 = "397829"
$ws.Range("G7").Value2 = "20190115"
$ws.Range("H7").Value2 = "SomeCT"
$ws.Range("I7").Value2 = "397829_CT1"
$ws.Range("L7").Value2 = "CT"
$ws.Range("O7").Value2 = "20190115"
$ws.Range("P7").Value2 = "CT1 abdomen"
$ws.Range("Q7").Value2 = "ROE^JANE"
$ws.Range("R7").Value2 = "397829"
$ws.Range("S7").Value2 = "SomeCT"
$ws.Range("T7").Value2 = "img-00000"
$ws.Range("A8").Value2 = "Scans"
$ws.Range("B8").Value2 = "tests/fixtures/basic/ROE^JANE-397829/20200623/SomeCT/img-00000.dcm"
$ws.Range("C8").Value2 = "img-00000.dcm"
$ws.Range("D8").Value2 = "Y"
$ws.Range("F8").Value2 = "397829"
$ws.Range("G8").Value2 = "20200623"
$ws.Range("H8").Value2 = "SomeCT"
$ws.Range("I8").Value2 = "397829_CT2"
$ws.Range("L8").Value2 = "CT"
$ws.Range("O8").Value2 = "20200623"
$ws.Range("P8").Value2 = "CT1 abdomen"
$ws.Range("Q8").Value2 = "ROE^JANE"
$ws.Range("R8").Value2 = "397829"
$ws.Range("S8").Value2 = "SomeCT"
$ws.Range("T8").Value2 = "img-00000"
$ws.Range("A9").Value2 = "Scans"
$ws.Range("B9").Value2 = "tests/fixtures/basic/ROE^JANE-397829/20210414/SomeCT/image-00000.dcm"
$ws.Range("C9").Value2 = "image-00000.dcm"
$ws.Range("D9").Value2 = "Y"
$ws.Range("F9").Value2 = "397829"
$ws.Range("G9").Value2 = "20210414"
$ws.Range("H9").Value2 = "SomeCT"
$ws.Range("I9").Value2 = "397829_CT3"
$ws.Range("L9").Value2 = "CT"
$ws.Range("O9").Value2 = "20210414"
$ws.Range("P9").Value2 = "CT1 abdomen"
$ws.Range("Q9").Value2 = "ROE^JANE"
$ws.Range("R9").Value2 = "397829"
$ws.Range("S9").Value2 = "SomeCT"
$ws.Range("T9").Value2 = "image-00000"
$ws.Range("A10").Value2 = "Scans"
$ws.Range("B10").Value2 = "tests/fixtures/basic/Smith^John-038945/20200303/X-Rays/img-00000.dcm"
$ws.Range("C10").Value2 = "img-00000.dcm"
$ws.Range("D10").Value2 = "Y"
$ws.Range("F10").Value2 = "038945"
$ws.Range("G10").Value2 = "20200303"
$ws.Range("H10").Value2 = "X-Rays"
$ws.Range("I10").Value2 = "038945_CT1"
$ws.Range("L10").Value2 = "CT"
$ws.Range("O10").Value2 = "20200303"
$ws.Range("P10").Value2 = "CT1 abdomen"
$ws.Range("Q10").Value2 = "Smith^John"
$ws.Range("R10").Value2 = "038945"
$ws.Range("S10").Value2 = "X-Rays"
$ws.Range("T10").Value2 = "img-00000"
$ws.Range("B11").Value2 = "tests/fixtures/basic/DOE^JOHN-002304/20200312HeadCT/Head CT/20200312-scan1.txt"
$ws.Range("C11").Value2 = "20200312-scan1.txt"
$ws.Range("D11").Value2 = "N"
$ws.Range("E11").Value2 = "unmatched"
$ws.Range("B12").Value2 = "tests/fixtures/basic/NomatchDir/no_match_file.txt"
$ws.Range("C12").Value2 = "no_match_file.txt"
$ws.Range("D12").Value2 = "N"
$ws.Range("E12").Value2 = "unmatched"
$ws.Range("B13").Value2 = "tests/fixtures/basic/ROE^JANE-397829/20190115/SomeCT/20190115-scan1.txt"
$ws.Range("C13").Value2 = "20190115-scan1.txt"
$ws.Range("D13").Value2 = "N"
$ws.Range("E13").Value2 = "unmatched"
$ws.Range("B14").Value2 = "tests/fixtures/basic/ROE^JANE-397829/20200623/SomeCT/20200623-scan1.txt"
$ws.Range("C14").Value2 = "20200623-scan1.txt"
$ws.Range("D14").Value2 = "N"
$ws.Range("E14").Value2 = "unmatched"
$ws.Range("B15").Value2 = "tests/fixtures/basic/ROE^JANE-397829/20210414/SomeCT/20210414-scan1.txt"
$ws.Range("C15").Value2 = "20210414-scan1.txt"
$ws.Range("D15").Value2 = "N"
$ws.Range("E15").Value2 = "unmatched"
$ws.Range("B16").Value2 = "tests/fixtures/basic/Smith^John-038945/20200303/X-Rays/20200303-scan1.txt"
$ws.Range("C16").Value2 = "20200303-scan1.txt"
$ws.Range("D16").Value2 = "N"
$ws.Range("E16").Value2 = "unmatched"
$ws.Range("B17").Value2 = "tests/fixtures/basic/Smith^John-038945/no_match_file.txt"
$ws.Range("C17").Value2 = "no_match_file.txt"
$ws.Range("D17").Value2 = "N"
$ws.Range("E17").Value2 = "unmatched"
# Remove the temporary text formatting again so the cells fall back to the
# workbook's default (unformatted) style, then fill in the numeric
# DICOM:SeriesNumber values (these must stay real numbers, not text).
$dataRange.ClearFormats()

$ws.Range("M2").Value2 = 6168
$ws.Range("M3").Value2 = 6168
$ws.Range("M4").Value2 = 6168
$ws.Range("M5").Value2 = 6168
$ws.Range("M6").Value2 = 6168
$ws.Range("M7").Value2 = 6168
$ws.Range("M8").Value2 = 6168
$ws.Range("M9").Value2 = 6168
$ws.Range("M10").Value2 = 6168
# Match the new active selection left by the refactored scan code.
$ws.Range("I2:I10").Select()
